# Super fast code with proper OOP
# Updates the ReqPow_AC / ReqPow_FC / ReqPow_Batt rows (B:U) with recalculated
# required-power values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ReqPow_AC), columns B:U -------------------------------------
$row2 = @(
    88.4545441393966, 176.9090882787932, 1101.96296801947, 1769.090882787932,
    1769.090882787932, 1590.364055129905, 1590.364055129905, 1484.992756861953,
    1484.992756861953, 176.9090882787932, 1769.090882787932, 1769.090882787932,
    1484.992756861953, 1484.992756861953, 1484.992756861953, 1484.992756861953,
    176.9090882787932, 176.9090882787932, 176.9090882787932, 88.4545441393966
)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, 2 + $i).Value = $row2[$i]
}

# --- Row 3 (ReqPow_FC), columns B:U -------------------------------------
$row3 = @(
    88.4545441393966, 176.9090882787932, 886.5192438262889, 1437.033505445532,
    1437.033505445532, 1645.013641234165, 1645.013641234165, 1484.992756861953,
    1484.992756861953, 176.9090882787932, 1437.033505445532, 1437.033505445532,
    1484.992756861953, 1484.992756861953, 1484.992756861953, 1484.992756861953,
    176.9090882787932, 176.9090882787932, 176.9090882787932, 88.4545441393966
)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, 2 + $i).Value = $row3[$i]
}

# --- Row 4 (ReqPow_Batt) -- only the non-zero columns changed ----------
$ws.Range("D4").Value = 215.443724193181
$ws.Range("E4").Value = 332.0573773423994
$ws.Range("F4").Value = 332.0573773423994
$ws.Range("G4").Value = -54.64958610426001
$ws.Range("H4").Value = -54.64958610426001
$ws.Range("L4").Value = 332.0573773423994
$ws.Range("M4").Value = 332.0573773423994
